$wb = $excel.ActiveWorkbook

# Sheet "TillattOrganisasjonsform": remove the "ESEK" organisasjonsform row (row 5)
$wsOrg = $wb.Worksheets.Item("TillattOrganisasjonsform")
$wsOrg.Rows.Item(5).Delete()

# Sheet "TillattOrganisasjonsformPosisjo": remove the 4 "ESEK" rows (rows 10-13)
$wsPos = $wb.Worksheets.Item("TillattOrganisasjonsformPosisjo")
$wsPos.Range("A10:A13").EntireRow.Delete()

# Fix up selections / active sheet to match the target state
$wsPos.Range("A16").Select()
$wsOrg.Activate()
$wsOrg.Range("B28").Select()
